$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update conscientiousness (column W) scores for rows 2-37 (row 24 has no data)
$ws.Range("W2").Value = 5
$ws.Range("W3").Value = 3
$ws.Range("W4").Value = 6
$ws.Range("W5").Value = 3
$ws.Range("W6").Value = 8
$ws.Range("W7").Value = 4
$ws.Range("W8").Value = 8
$ws.Range("W9").Value = 6
$ws.Range("W10").Value = 2
$ws.Range("W11").Value = 6
$ws.Range("W12").Value = 3
$ws.Range("W13").Value = 8
$ws.Range("W14").Value = 8
$ws.Range("W15").Value = 6
$ws.Range("W16").Value = 8
$ws.Range("W17").Value = 3
$ws.Range("W18").Value = 5
$ws.Range("W19").Value = 5
$ws.Range("W20").Value = 8
$ws.Range("W21").Value = 8
$ws.Range("W22").Value = 3
$ws.Range("W23").Value = 5
$ws.Range("W25").Value = 3
$ws.Range("W26").Value = 7
$ws.Range("W27").Value = 4
$ws.Range("W28").Value = 6
$ws.Range("W29").Value = 8
$ws.Range("W30").Value = 8
$ws.Range("W31").Value = 9
$ws.Range("W32").Value = 7
$ws.Range("W33").Value = 5
$ws.Range("W34").Value = 6
$ws.Range("W35").Value = 10
$ws.Range("W36").Value = 7
$ws.Range("W37").Value = 4

# Re-write coffee (column AR) values as numbers instead of booleans for rows 29-50
$ws.Range("AR29").Value = 0
$ws.Range("AR30").Value = 1
$ws.Range("AR31").Value = 0
$ws.Range("AR32").Value = 0
$ws.Range("AR33").Value = 0
$ws.Range("AR34").Value = 0
$ws.Range("AR35").Value = 0
$ws.Range("AR36").Value = 0
$ws.Range("AR37").Value = 1
$ws.Range("AR38").Value = 0
$ws.Range("AR39").Value = 1
$ws.Range("AR40").Value = 0
$ws.Range("AR41").Value = 0
$ws.Range("AR42").Value = 0
$ws.Range("AR43").Value = 0
$ws.Range("AR44").Value = 0
$ws.Range("AR45").Value = 1
$ws.Range("AR46").Value = 0
$ws.Range("AR47").Value = 0
$ws.Range("AR48").Value = 1
$ws.Range("AR49").Value = 1
$ws.Range("AR50").Value = 0
